$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.0674
$ws.Range("E2").Value = 0.138
$ws.Range("F2").Value = 0.0673
$ws.Range("G2").Value = 0.4367856961446639
$ws.Range("H2").Value = 0.2671305937928582
$ws.Range("I2").Value = 0.2446732365725875
$ws.Range("J2").Value = 0.2446732365725875
$ws.Range("K2").Value = 1019.6
$ws.Range("L2").Value = 0.2589526083202113
$ws.Range("M2").Value = 468.3
$ws.Range("N2").Value = 0.02046452684248476
$ws.Range("O2").Value = 0.4592977638289525
$ws.Range("P2").Value = 442.3
$ws.Range("Q2").Value = 0.01932833701138375
$ws.Range("R2").Value = 0.4337975676735975
$ws.Range("S2").Value = 26
$ws.Range("T2").Value = 0.05551996583386718
$ws.Range("U2").Value = 1223.5
$ws.Range("V2").Value = 0.0534664714750803
$ws.Range("W2").Value = 0.2327375653404552
$ws.Range("X2").Value = 0.0634871818159634
$ws.Range("Y2").Value = 0.1692503835244918
$ws.Range("Z2").Value = 1.133678731770138
$ws.Range("AA2").Value = 0.2773808445357058
$ws.Range("AB2").Value = 0.0633462619521882
$ws.Range("AC2").Value = 0.2140345825835176
$ws.Range("AD2").Value = 0
$ws.Range("AE2").Value = 68.61799159547031
$ws.Range("AF2").Value = 68.61799159547031
$ws.Range("AG2").Value = -1154.88200840453
$ws.Range("AH2").Value = 0.002989614798102581
$ws.Range("AI2").Value = 0.01326485200557116
$ws.Range("AJ2").Value = -0.05315027439164482
$ws.Range("AK2").Value = -0.2924182780506312
$ws.Range("AM2").Value = -43.3
$ws.Range("AN2").Value = 0
$ws.Range("AP2").Value = -1.05391678080355
$ws.Range("AQ2").Value = -22.15935334872979

# Row 3
$ws.Range("D3").Value = 0.0674
$ws.Range("E3").Value = 0.138
$ws.Range("F3").Value = 0.0673
$ws.Range("G3").Value = 0.4367856961446639
$ws.Range("H3").Value = 0.2671305937928582
$ws.Range("I3").Value = 0.2446732365725875
$ws.Range("J3").Value = 0.2446732365725875
$ws.Range("K3").Value = 1019.6
$ws.Range("L3").Value = 0.2589526083202113
$ws.Range("M3").Value = 468.3
$ws.Range("N3").Value = 0.02046452684248476
$ws.Range("O3").Value = 0.4592977638289525
$ws.Range("P3").Value = 442.3
$ws.Range("Q3").Value = 0.01932833701138375
$ws.Range("R3").Value = 0.4337975676735975
$ws.Range("S3").Value = 26
$ws.Range("T3").Value = 0.05551996583386718
$ws.Range("U3").Value = 1223.5
$ws.Range("V3").Value = 0.0534664714750803
$ws.Range("W3").Value = 0.2327375653404552
$ws.Range("X3").Value = 0.0634871818159634
$ws.Range("Y3").Value = 0.1692503835244918
$ws.Range("Z3").Value = 1.133678731770138
$ws.Range("AA3").Value = 0.2773808445357058
$ws.Range("AB3").Value = 0.0633462619521882
$ws.Range("AC3").Value = 0.2140345825835176
$ws.Range("AD3").Value = 0
$ws.Range("AE3").Value = 68.61799159547031
$ws.Range("AF3").Value = 68.61799159547031
$ws.Range("AG3").Value = -1154.88200840453
$ws.Range("AH3").Value = 0.002989614798102581
$ws.Range("AI3").Value = 0.01326485200557116
$ws.Range("AJ3").Value = -0.05315027439164482
$ws.Range("AK3").Value = -0.2924182780506312
$ws.Range("AM3").Value = -43.3
$ws.Range("AN3").Value = 0
$ws.Range("AP3").Value = -1.05391678080355
$ws.Range("AQ3").Value = -22.15935334872979

